# Replace the hardcoded dynamic UUID-based path parameters with
# placeholders that are resolved from system properties at test time.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D2 / D5: "/follow/user/<user1>/following/<user2>"
$ws.Range("D2").Value = "/follow/user/(SYS_USER1)/following/(SYS_USER2)"
$ws.Range("D5").Value = "/follow/user/(SYS_USER1)/following/(SYS_USER2)"

# D3: "/follow/user/<user1>/following"
$ws.Range("D3").Value = "/follow/user/(SYS_USER1)/following"

# D4: "/follow/user/<user2>/followers"
$ws.Range("D4").Value = "/follow/user/(SYS_USER2)/followers"

# D2 previously carried the "Hyperlink" cell style (blue/underlined) even
# though it was never a real hyperlink. Since it's no longer a literal
# URL, drop that formatting back to the default "Normal" style.
$ws.Range("D2").Style = "Normal"

# The workbook no longer has any cell using the "Hyperlink" style, so
# remove the now-unused named style (and its backing font/format) from
# the workbook's style table entirely.
try {
    $wb.Styles("Hyperlink").Delete()
} catch {
    # Already absent (e.g. script re-run) - nothing to clean up.
}

# Reset the view: clear the scrolled-right "topLeftCell" and move the
# selection back to A2 instead of the old L2:L5 block.
$ws.Range("A2").Select()
